$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 521.8570999999999
$ws.Range("I58").Value = 593.8333
$ws.Range("J58").Value = 90
$ws.Range("K58").Value = 1781.4999
$ws.Range("L58").Value = 270
$ws.Range("M58").Value = -1631.4999
$ws.Range("N58").Value = -570

$ws.Range("H86").Value = 2608.0908
$ws.Range("I86").Value = 1693
$ws.Range("J86").Value = 4209.5
$ws.Range("K86").Value = 1693
$ws.Range("L86").Value = 4209.5
$ws.Range("M86").Value = -570
$ws.Range("N86").Value = -6455.5

$ws.Range("H87").Value = 179998.33
$ws.Range("I87").Value = 40000
$ws.Range("J87").Value = 249997.5
$ws.Range("K87").Value = 40000
$ws.Range("L87").Value = 249997.5
$ws.Range("M87").Value = -38752
$ws.Range("N87").Value = -252493.5

$ws.Range("H89").Value = 2608.0908
$ws.Range("I89").Value = 1693
$ws.Range("J89").Value = 4209.5
$ws.Range("K89").Value = 8465
$ws.Range("L89").Value = 21047.5
$ws.Range("M89").Value = -2849
$ws.Range("N89").Value = -32279.5

$ws.Range("H90").Value = 179998.33
$ws.Range("I90").Value = 40000
$ws.Range("J90").Value = 249997.5
$ws.Range("K90").Value = 120000
$ws.Range("L90").Value = 749992.5
$ws.Range("M90").Value = -113760
$ws.Range("N90").Value = -762472.5

$ws.Range("H112").Value = 3102.06
$ws.Range("I112").Value = 1747.5
$ws.Range("J112").Value = 3158.5
$ws.Range("K112").Value = 5242.5
$ws.Range("L112").Value = 9475.5
$ws.Range("M112").Value = -4134.5
$ws.Range("N112").Value = -11691.5

$ws.Range("H137").Value = 1333.3448
$ws.Range("I137").Value = 1276.25
$ws.Range("J137").Value = 1607.4
$ws.Range("K137").Value = 3828.75
$ws.Range("L137").Value = 4822.200000000001
$ws.Range("M137").Value = -1278.75
$ws.Range("N137").Value = -9922.200000000001

$ws.Range("H138").Value = 5480.5557
$ws.Range("I138").Value = 1899
$ws.Range("J138").Value = 6988.579
$ws.Range("K138").Value = 5697
$ws.Range("L138").Value = 20965.737
$ws.Range("M138").Value = -557
$ws.Range("N138").Value = -31245.737

$ws.Range("H141").Value = 2335
$ws.Range("I141").Value = 2311.5715
$ws.Range("J141").Value = 2499
$ws.Range("K141").Value = 6934.7145
$ws.Range("L141").Value = 7497
$ws.Range("M141").Value = -1754.7145
$ws.Range("N141").Value = -17857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 8619.571
$ws.Range("I16").Value = 6723
$ws.Range("J16").Value = 19999
$ws.Range("K16").Value = 6723
$ws.Range("L16").Value = 19999
$ws.Range("M16").Value = -6436
$ws.Range("N16").Value = -20573

$ws.Range("H61").Value = 4497.6343
$ws.Range("I61").Value = 3279.8235
$ws.Range("J61").Value = 10412.714
$ws.Range("K61").Value = 3279.8235
$ws.Range("L61").Value = 10412.714
$ws.Range("M61").Value = -3067.8235
$ws.Range("N61").Value = -10836.714

$ws.Range("H74").Value = 2131.8064
$ws.Range("I74").Value = 2183.577
$ws.Range("J74").Value = 1862.6
$ws.Range("K74").Value = 2183.577
$ws.Range("L74").Value = 1862.6
$ws.Range("M74").Value = -1309.577
$ws.Range("N74").Value = -3610.6

$ws.Range("H77").Value = 2131.8064
$ws.Range("I77").Value = 2183.577
$ws.Range("J77").Value = 1862.6
$ws.Range("K77").Value = 10917.885
$ws.Range("L77").Value = 9313
$ws.Range("M77").Value = -6549.885000000002
$ws.Range("N77").Value = -18049

$ws.Range("H136").Value = 4497.6343
$ws.Range("I136").Value = 3279.8235
$ws.Range("J136").Value = 10412.714
$ws.Range("K136").Value = 9839.470499999999
$ws.Range("L136").Value = 31238.142
$ws.Range("M136").Value = -7289.470499999999
$ws.Range("N136").Value = -36338.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7040.2964
$ws.Range("I31").Value = 13398.692
$ws.Range("J31").Value = 5024.2197
$ws.Range("K31").Value = 13398.692
$ws.Range("L31").Value = 5024.2197
$ws.Range("M31").Value = -13103.692
$ws.Range("N31").Value = -5614.2197

$ws.Range("H34").Value = 7040.2964
$ws.Range("I34").Value = 13398.692
$ws.Range("J34").Value = 5024.2197
$ws.Range("K34").Value = 13398.692
$ws.Range("L34").Value = 5024.2197
$ws.Range("M34").Value = -13196.692
$ws.Range("N34").Value = -5428.2197

$ws.Range("H58").Value = 2761.6667
$ws.Range("I58").Value = 2705.5
$ws.Range("J58").Value = 2874
$ws.Range("K58").Value = 2705.5
$ws.Range("L58").Value = 2874
$ws.Range("M58").Value = -2502.5
$ws.Range("N58").Value = -3280

$ws.Range("H107").Value = 22727962
$ws.Range("I107").Value = 35714988
$ws.Range("J107").Value = 669.625
$ws.Range("K107").Value = 35714988
$ws.Range("L107").Value = 669.625
$ws.Range("M107").Value = -35713068
$ws.Range("N107").Value = -4509.625

$ws.Range("H122").Value = 139769.36
$ws.Range("I122").Value = 218008.64
$ws.Range("J122").Value = 2850.625
$ws.Range("K122").Value = 654025.92
$ws.Range("L122").Value = 8551.875
$ws.Range("M122").Value = -651575.92
$ws.Range("N122").Value = -13451.875

$ws.Range("H132").Value = 1212.7142
$ws.Range("I132").Value = 1148.8628
$ws.Range("J132").Value = 1864
$ws.Range("K132").Value = 3446.588400000001
$ws.Range("L132").Value = 5592
$ws.Range("M132").Value = -916.5884000000005
$ws.Range("N132").Value = -10652

$ws.Range("H134").Value = 962.4138
$ws.Range("I134").Value = 900.8182
$ws.Range("J134").Value = 1156
$ws.Range("K134").Value = 2702.4546
$ws.Range("L134").Value = 3468
$ws.Range("M134").Value = -167.4546
$ws.Range("N134").Value = -8538

$ws.Range("H136").Value = 2761.6667
$ws.Range("I136").Value = 2705.5
$ws.Range("J136").Value = 2874
$ws.Range("K136").Value = 8116.5
$ws.Range("L136").Value = 8622
$ws.Range("M136").Value = -5566.5
$ws.Range("N136").Value = -13722

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3494.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3494.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10483.5
$ws.Range("N81").Value = -12729.5

$ws.Range("H84").Value = 3494.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3494.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 31450.5
$ws.Range("N84").Value = -42682.5

$ws.Range("H117").Value = 5558.091
$ws.Range("I117").Value = 3263.75
$ws.Range("J117").Value = 6869.143
$ws.Range("K117").Value = 9791.25
$ws.Range("L117").Value = 20607.429
$ws.Range("M117").Value = -6349.25
$ws.Range("N117").Value = -27491.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5078.476
$ws.Range("I80").Value = 3978.5715
$ws.Range("J80").Value = 7278.2856
$ws.Range("K80").Value = 3978.5715
$ws.Range("L80").Value = 7278.2856
$ws.Range("M80").Value = -2980.5715
$ws.Range("N80").Value = -9274.285599999999

$ws.Range("H83").Value = 5078.476
$ws.Range("I83").Value = 3978.5715
$ws.Range("J83").Value = 7278.2856
$ws.Range("K83").Value = 19892.8575
$ws.Range("L83").Value = 36391.428
$ws.Range("M83").Value = -14900.8575
$ws.Range("N83").Value = -46375.428

$ws.Range("H119").Value = 89000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 89000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 89000
$ws.Range("N119").Value = -98676

$ws.Range("H132").Value = 5320.0835
$ws.Range("I132").Value = 5508.7144
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 16526.1432
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -13996.1432
$ws.Range("N132").Value = -17059.0001

$ws.Range("H136").Value = 30522.125
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 30522.125
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 91566.375
$ws.Range("N136").Value = -96666.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5900.364
$ws.Range("I22").Value = 4780.6665
$ws.Range("J22").Value = 10939
$ws.Range("K22").Value = 4780.6665
$ws.Range("L22").Value = 10939
$ws.Range("M22").Value = -4485.6665
$ws.Range("N22").Value = -11529

$ws.Range("H25").Value = 9000
$ws.Range("I25").Value = 9000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -8770

$ws.Range("H27").Value = 5900.364
$ws.Range("I27").Value = 4780.6665
$ws.Range("J27").Value = 10939
$ws.Range("K27").Value = 4780.6665
$ws.Range("L27").Value = 10939
$ws.Range("M27").Value = -4673.6665
$ws.Range("N27").Value = -11153

$ws.Range("H40").Value = 2564.9524
$ws.Range("I40").Value = 2144.375
$ws.Range("J40").Value = 3910.8
$ws.Range("K40").Value = 2144.375
$ws.Range("L40").Value = 3910.8
$ws.Range("M40").Value = -2008.375
$ws.Range("N40").Value = -4182.8

$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 20000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -20856

$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 20000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -22964

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 67250
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 67250
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 67250
$ws.Range("N56").Value = -68678

$ws.Range("H107").Value = 339.2857
$ws.Range("I107").Value = 346
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 1038
$ws.Range("L107").Value = 897
$ws.Range("M107").Value = 882
$ws.Range("N107").Value = -4737

$ws.Range("H132").Value = 4678.225
$ws.Range("I132").Value = 2639.6333
$ws.Range("J132").Value = 10794
$ws.Range("K132").Value = 7918.8999
$ws.Range("L132").Value = 32382
$ws.Range("M132").Value = -5388.8999
$ws.Range("N132").Value = -37442

$ws.Range("H136").Value = 3438
$ws.Range("I136").Value = 3746
$ws.Range("J136").Value = 2668
$ws.Range("K136").Value = 11238
$ws.Range("L136").Value = 8004
$ws.Range("M136").Value = -8688
$ws.Range("N136").Value = -13104
